$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.543.12'
$ws.Range("E2").Value = '  -0.14%  '
$ws.Range("D3").Value = '3.046.02'
$ws.Range("E3").Value = '  -0.98%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = "'554.16"
$ws.Range("E5").Value = '  +0.58%  '
$ws.Range("D6").Value = "'141.80"
$ws.Range("E6").Value = '  -1.21%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").Value = '3.047.80'
$ws.Range("E8").Value = '  -0.78%  '
$ws.Range("D9").Value = "'0.501"
$ws.Range("E9").Value = '  -0.49%  '
$ws.Range("D10").Value = "'0.153"
$ws.Range("E10").Value = '  +0.06%  '
$ws.Range("D11").Value = "'6.06"
$ws.Range("E11").Value = '  -7.14%  '
$ws.Range("D12").Value = "'0.467"
$ws.Range("E12").Value = '  +1.71%  '
$ws.Range("D13").Value = "'0.0000226"
$ws.Range("E13").Value = '  -1.04%  '
$ws.Range("D14").Value = "'34.58"
$ws.Range("E14").Value = '  -1.34%  '
$ws.Range("D15").Value = '3.566.85'
$ws.Range("E15").Value = '  -0.01%  '
$ws.Range("D16").Value = '63.591.51'
$ws.Range("E16").Value = '  +0.12%  '
$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").Value = "'0.110"
$ws.Range("E17").Value = '  +0.42%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.039.03'
$ws.Range("E18").Value = '  -1.12%  '
$ws.Range("D19").Value = "'6.65"
$ws.Range("E19").Value = '  -2.39%  '
$ws.Range("D20").Value = "'475.21"
$ws.Range("E20").Value = '  -2.06%  '
$ws.Range("D21").Value = "'13.88"
$ws.Range("E21").Value = '  -0.35%  '
$ws.Range("D22").Value = "'0.669"
$ws.Range("E22").Value = '  -1.57%  '
$ws.Range("D23").Value = "'7.45"
$ws.Range("E23").Value = '  +1.70%  '
$ws.Range("D24").Value = "'13.95"
$ws.Range("E24").Value = '  +9.00%  '
$ws.Range("D25").Value = "'80.72"
$ws.Range("E25").Value = '  -0.67%  '
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("D27").Value = "'2.76"
$ws.Range("E27").Value = '  -0.80%  '
$ws.Range("D28").Value = "'7.88"
$ws.Range("E28").Value = '  -1.46%  '
$ws.Range("D29").Value = "'2.03"
$ws.Range("E29").Value = '  +0.30%  '
$ws.Range("D30").Value = "'0.996"
$ws.Range("E30").Value = '  -0.24%  '
$ws.Range("D31").Value = "'26.02"
$ws.Range("E31").Value = '  -1.01%  '
$ws.Range("D32").Value = "'1.13"
$ws.Range("E32").Value = '  -3.33%  '
$ws.Range("D33").Value = "'2.43"
$ws.Range("E33").Value = '  -0.95%  '
$ws.Range("D34").Value = "'5.56"
$ws.Range("E34").Value = '  -2.86%  '
$ws.Range("D35").Value = "'6.12"
$ws.Range("E35").Value = '  +1.39%  '
$ws.Range("D36").Value = "'54.98"
$ws.Range("E36").Value = '  -1.32%  '
$ws.Range("D37").Value = "'0.0403"
$ws.Range("E37").Value = '  +0.79%  '
$ws.Range("D38").Value = "'2.87"
$ws.Range("E38").Value = '  +10.92%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = "'0.0801"
$ws.Range("E39").Value = '  -2.93%  '
$ws.Range("B40").Value = 'Bittensor'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D40").Value = "'432.65"
$ws.Range("E40").Value = '  -7.87%  '
$ws.Range("D41").Value = '2.935.95'
$ws.Range("E41").Value = '  -3.26%  '
$ws.Range("D42").Value = "'8.11"
$ws.Range("E42").Value = '  -1.72%  '
$ws.Range("D43").Value = "'0.112"
$ws.Range("E43").Value = '  -5.88%  '
$ws.Range("D44").Value = "'27.93"
$ws.Range("E44").Value = '  +0.42%  '
$ws.Range("D45").Value = "'0.257"
$ws.Range("E45").Value = '  +0.40%  '
$ws.Range("E46").Value = '  +0.00%  '
$ws.Range("D47").Value = "'2.10"
$ws.Range("E47").Value = '  +1.97%  '
$ws.Range("D48").Value = "'0.111"
$ws.Range("E48").Value = '  +0.23%  '
$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").Value = "'116.96"
$ws.Range("E49").Value = '  +0.19%  '
$ws.Range("B50").Value = 'PEPE'
$ws.Range("C50").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D50").Value = '0.0₃0508'
$ws.Range("E50").Value = '  -1.04%  '
$ws.Range("D51").Value = "'2.05"
$ws.Range("E51").Value = '  -2.09%  '
